$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @("en", "fr", "status"),
  @(".acorn data successfully generated!", "Le fichier .acorn a bien été généré", ""),
  @(".acorn file saved on server.", "Fichier .acorn sauvegardé sur le serveur.", ""),
  @(".acorn not saved.", "Fichier .acorn n'est pas sauvegardé.", ""),
  @("(1/4) Download Clinical data", "(1/4) Télécharger les Données Clinique", ""),
  @("(2/4) Provide Lab data", "(2/4) Fournir des données de laboratoire", ""),
  @("(3/4) Combine Clinical and Lab data", "(3/4) Combiner les données clinique et de laboratoire", ""),
  @("(4/4) Save .acorn file", "(4/4) Sauvegarder un ficher .acorn", ""),
  @("(Optional) Comments:", "(Optionnel) Commentaires:", ""),
  @("(To log out, close the app.)", "(Pour vous déconnecter, fermer l'app.)", ""),
  @("ACORN Participating Countries", "Pays Participants au Projet ACORN", ""),
  @("All 'orgname' are provided.", "Tous les 'orgname' sont fournis.", ""),
  @("All 'patid' are provided.", "Tous les 'patid' sont fournis.", ""),
  @("All 'specdate' are provided.", "Tous les 'specdate' sont fournis.", ""),
  @("All 'specdate' are today or before today.", "Tous les 'specdate' sont égales ou antérieures à aujourd'hui.", ""),
  @("All 'specgroup' are provided.", "Tous les 'specgroup' sont fournis.", ""),
  @("All 'specid' are provided.", "Tous les 'specid' sont fournis.", ""),
  @("All dates of enrolment for HAI patients have a matching date in the HAI survey dataset", "Toutes les dates d'enrôlements de patients avec une HAI se trouvent dans le jeu de données d'études de HAI.", ""),
  @("All Other Organisms", "Tous les Autres Organismes", ""),
  @("All valid records have an ACORN ID.", "Tous les enregistrements valides ont une ID ACORN.", ""),
  @("AMR", "AMR", ""),
  @("and generate enrolment log.", "et générer un fichier de suivi des recrutements.", ""),
  @("Attempting to connect.", "Tentative de connection.", ""),
  @("Blood culture collected within 24 hours of admission (CAI) / symptom onset (HAI)", "Hémoculture collectée dans les 24 heures d'admission (CAI) / au début des symptômes (HAI)", ""),
  @("Blood Culture Contaminants", "Contaminants d'Hémoculture", ""),
  @("Bloodstream Infection (BSI)", "Infection du Sang (BSI)", ""),
  @("Calculated age is consistent with 'Age Category'", "L'âge calculé est cohérent avec la catégorie d'âge.", ""),
  @("Calculated age isn't always consistent with 'Age Category'", "L'âge calculé n'est pas toujours cohérent avec la catégorie d'âge.", ""),
  @("Cancel", "Annuler", ""),
  @("Care should be taken when interpreting rates and AMR profiles where there are small numbers of cases or bacterial isolates: point estimates may be unreliable.", "L'interprétation des taux et profils AMR doit prêter à caution lorsqu'en présence d'un petit nombre de cas ou d'isolats : les estimations sont alors peu fiable.", ""),
  @("Clinical and day-28 outcomes are consistent.", "Les résultats cliniques et du jour 28 sont cohérents.", ""),
  @("Clinical and day-28 outcomes aren't consistent for some dead patients.", "Pour certains patients décédés, les résultats cliniques et du jour 28 ne sont pas cohérents.", ""),
  @("Clinical Outcome", "Résultat Clinique", ""),
  @("Clinical Outcome Status:", "Statut des résultats cliniques:", ""),
  @("Co-resistances", "Co-résistances", ""),
  @("Combine Susceptible + Intermediate", "Combiner Susceptible + Intermédiaire", ""),
  @("Consider saving .acorn file on the cloud for additional security.", "Pour plus de sécurité, envisagez une sauvegarde des données sur un serveur.", ""),
  @("Contains names of organisms before and after mapping.", "Contient les noms des organismes avant et après cartorgraphie.", ""),
  @("Couldn't connect to server. Please check internet access.", "Il n'est pas possible de se connecter au server. Veuillez vérifier la connection à internet", ""),
  @("Critical errors with clinical data.", "Erreurs critiques avec les données cliniques.", ""),
  @("Culture results per specimen type", "Résultats de culture par type de spécimen", ""),
  @("Data Management", "Gestion des Données", ""),
  @("Date of Enrolment", "Date d'Enrôlement", ""),
  @("Day 28", "Jour 28", ""),
  @("Day 28 Status:", "Statut au jour 28:", ""),
  @("Diagnosis at Enrolment", "Diagnostics à l'Enrôlement", ""),
  @("Dismiss", "Annuler", ""),
  @("Distribution of Enrolments", "Distribution des Enrôlements", ""),
  @("Download Enrolment Log (.xlsx)", "Télécharger un fichier de suivi des recrutements (.xlsx)", ""),
  @("Download Lab Log (.xlsx)", "Télécharger un log du traitement des donnés de lab (.xlsx)", ""),
  @("Empiric Antibiotics Prescribed", "Antibiotiques Empiriques Prescris", ""),
  @("Enrolments", "Enrôlements", ""),
  @("Enrolments by (type of) Ward", "Enrôlements par (type de) service", ""),
  @("Enrolments with Blood Culture", "Enrôlements avec hémoculture", ""),
  @("Error in combining clinical and lab data.", "Erreur produite en combinant les données clinique et de laboratoire.", ""),
  @("Every D28 record (F04) matches exactly one patient enrolment (F01).", "Chaque enregistrement D28 (F04) correspond avec exactement un enrôlement de patient (F01).", ""),
  @("Every hospital outcome record (F03) has a matching infection episode (F02).", "A tout enregistrement de résultats hospitaliers (F03) correspond un épisode infectieux (F02).", ""),
  @("Every hospital outcome record (F03) has a matching patient enrolment (F01).", "A tout enregistrement de résultats hospitaliers (F03) correspond un enrôlement de patient (F01)", ""),
  @("Every infection episode record (F02) has a matching patient enrolment (F01).", "A tout enregistrement d'épisode infectieux (F02) correspond un enrôlement de patient (F01)", ""),
  @("File name:", "Nom de fichier:", ""),
  @("First sheet is the log of all enrolments retrived from REDCap (as per adjacent table). The second sheet is a listing of all flagged elements.", "Le premier onglet comprend tous les recrutements téléchargés depuis REDCap (voir table ci-contre). Le second onglet liste tous les éléments potentiellement à corriger.", ""),
  @("Follow-up", "Suivi", ""),
  @("from cultures that have growth", "de cultures avec croissance", ""),
  @("Generate .acorn file", "Générer un fichier .acorn", ""),
  @("Generate and load .acorn from clinical and lab data", "Générer un fichier .acorn depuis les données cliniques et de lab", ""),
  @("Generating .acorn", "Création du fichier .acorn", ""),
  @("Get data from REDCap", "Obtenir des données depuis REDCap", ""),
  @("Get the latest production release", "Obtenir la dernière version stable", ""),
  @("Growth / No Growth", "Croissance / Pas de croissance", ""),
  @("HAI point prevalence by ", "TBT", "new"),
  @("HAI Prevalence", "Prévalence des HAI", ""),
  @("Horizontal bars show the size of a set of SR results while vertical bars show the number of resistant isolates for the corresponding antibiotic.", "Les barres horizontales indiquent la taille d'un ensemble de résultats SR tandis que les barres verticales indiquent le nombre d'isolats résistants pour l'antibiotique correspondant.", ""),
  @("Info on loaded .acorn", "Information sur le fichier .acorn chargé.", ""),
  @("Initial & Final Surveillance Diagnosis", "Diagnostic de surveillance initial et final", ""),
  @("Isolates", "Isolats", ""),
  @("Issue detected with REDCap data. Please report to ACORN data managers. Until resolution, only existing .acorn files can be used.", "Problème détecté avec les données REDCap. Merci de contacter l'équipe ACORN. Jusqu'à résolution, seuls les fichiers .acorn existants peuvent être utilisés.", ""),
  @("It might take a couple of minutes. This window will close on completion.", "Cela peut prendre plusieurs minutes. Cette fenêtre se fermera automatiquement.", ""),
  @("Lab data successfully processed!", "Traitement des données de lab achevé avec succès !", ""),
  @("Lab data successfully provided.", "Les données de laboratoire ont été fournies avec succès.", ""),
  @("Lab dataset contains the minimal columns.", "Le jeu de données de lab contient les colonnes minimales requises.", ""),
  @("Lab dataset does not contains the minimal columns.", "Le jeu de données de lab ne contient pas les colonnes minimales requises.", ""),
  @("Language", "Langue", ""),
  @("Load .acorn", "Charger le fichier .acorn", ""),
  @("Load .acorn from cloud", "Charger le fichier .acorn depuis le nuage", ""),
  @("Load .acorn from local file", "Charger le fichier .acorn localement", ""),
  @("Load selected .acorn", "Charger le fichier .acorn sélectionné", ""),
  @("Loading data.", "Chargement des données.", ""),
  @("Log in", "Se connecter", ""),
  @("Microbiology", "Microbiologie", ""),
  @("Most frequent 10 organisms in the plot and complete listing in the table. Contaminants are in red.", "Les 10 organismes les plus nombreux sont dans le graphique et la liste complète dans la tables. Les contaminants sont en rouge.", ""),
  @("No .acorn data loaded.", "Aucun fichier .acorn n'a été chargé.", ""),
  @("No Blood Culture", "Pas d'Hématoculture", ""),
  @("Not connected to internet.", "Pas de connection à internet détecté.", ""),
  @("Number of specimens per specimen type", "Nombre de specimens par type de specimen", ""),
  @("Occupancy rate per type of ward per month", "Taux d'occupation par type de service et par mois", ""),
  @("of blood cultures grew a potential contaminant.", "des hémocultures ont générées un contaminant potentiel.", ""),
  @("of cultures have growth.", "des cultures avec croissance.", ""),
  @("of enrolments with blood culture.", "des enrôlements avec hémoculture.", ""),
  @("of Target Pathogens", "des Pathogènes Cibles", ""),
  @("Only isolates that have been tested against all of the drugs are included in the upset plot.", "Seuls les isolats qui ont été testés contre tous les médicaments sont inclus dans le graphique.", ""),
  @("Overview", "Aperçu", ""),
  @("Password", "Mot de passe", ""),
  @("Patient Age Distribution", "Distribution des ages des patients", ""),
  @("Patient Comorbidities", "Comorbiditées", ""),
  @("Patient enrolments", "Enrôlements", ""),
  @("Patients Transferred", "Patients Transferrés", ""),
  @("Please log in", "Connection", ""),
  @("Processing lab data.", "Traitement des données de lab.", ""),
  @("Reading lab data.", "Lecture des données de lab.", ""),
  @("Remove 'Not Cultured' specimens", "Supprimer les spécimens 'Not Cultured'", ""),
  @("Remove blood culture contaminants from the following visualizations", "Élimine les contaminants d'hémoculture des visualisations ci-dessous", ""),
  @("Reset Enrolments Filters", "Réinitialiser les Filtres sur les Enrôlements", ""),
  @("Resistance to 3rd gen. Cephalosporins Over Time", "Evolution de la Résistance aux 3rd gen. Cephalosporins", ""),
  @("Resistance to Carbapenems Over Time", "Evolution de la Résistance aux Carbapenems", ""),
  @("Resistance to Fluoroquinolones Over Time", "Evolution de la Résistance aux Fluoroquinolones", ""),
  @("Resistance to Oxacillin Over Time", "Evolution de la Résistance aux Oxacillin", ""),
  @("Resistance to Penicillin G - meningitis Over Time", "Evolution de la Résistance aux Penicillin G - meningitis", ""),
  @("Resistance to Penicillin G Over Time", "Evolution de la Résistance aux Penicillin G", ""),
  @("Retriving data from REDCap server.", "Récupération des données depuis le serveur REDcap.", ""),
  @("Save .acorn file", "Sauvegarder un ficher .acorn", ""),
  @("Save acorn data", "Sauvegarder des données acorn", ""),
  @("Save on Server", "Sauvegarder sur serveur", ""),
  @("See Breakdown by Ward", "Montrer la Répartition par Service", ""),
  @("See by Week", "Montrer par Semaine", ""),
  @("Show antibiotics combinations", "Montrer les combinaisons d'antibiotiques.", ""),
  @("Show comorbidities combinations", "Montrer les combinaisons de comorbiditées", ""),
  @("SIR Evaluation", "Evaluation SIR", ""),
  @("Some D28 records (F04) don't have a matching patient enrolment (F01).", "Certains enregistrements D28 (F04) n'ont pas un enrôlement de patient (F01) correspondant.", ""),
  @("Some dates of enrolment for HAI patients do have a matching date in the HAI survey dataset", "Certaines dates d'enrôlements de patients avec une HAI ne se trouvent pas dans le jeu de données d'études de HAI.", ""),
  @("Some hospital outcome records (F03) don't have a matching infection episode (F02). These records have been removed.", "Certains enregistrements de résultats hospitaliers (F03) ne correspondent à aucun un épisode infectieux (F02).", ""),
  @("Some hospital outcome records (F03) don't have a matching patient enrolment (F01).", "Certains enregistrements de résultats hospitaliers (F03) ne correspondent à aucun enrôlement de patient (F01)", ""),
  @("Some infection episode records (F02) don't have a matching patient enrolment (F01). These records have been removed.", "Certains enregistrements d'épisode infectieux (F02) ne correspondent à aucun enrôlement de patient (F01). Ces enregistrements ont été supprimés.", ""),
  @("Some records with a missing ACORN ID. These records have been removed.", "Certains enregistrements ont un identifiant ACORN manquant. Ces enregistrements ont été supprimés.", ""),
  @("Specimen Types", "Types de Specimen", ""),
  @("Specimens", "Spécimens", ""),
  @("Specimens Collected", "Spécimens Recueillis", ""),
  @("specimens per enrolment", "spécimens par enrôlement", ""),
  @("Successfully combined clinical and lab data into .acorn file", "Les données cliniques et de laboratoire ont bien été combinées dans un fichier .acorn", ""),
  @("Successfully loaded data.", "Chargement des données réalisé avec succès.", ""),
  @("Successfully logged in.", "Connection réussie.", ""),
  @("Successfully saved .acorn file in the cloud. You can now explore acorn data.", "Fichier .acorn sauvegardé sur le serveur. Vous pouvez maintenant explorer les données acorn.", ""),
  @("Successfully saved .acorn file locally.", "Fichier .acorn sauvegardé locallement.", ""),
  @("Supply first valid clinical and lab data.", "Fournir en premier lieu des données clinique et de lab valides.", ""),
  @("Susceptible & Intermediate are always combined in this visualisation of co-resistances.", "Susceptible & Intermédiaire sont toujours combinés dans cette visualisation des co-résistances.", ""),
  @("The 10 most common initial-final diagnosis combinations:", "Les 10 combinaisons de diagnostic initial-final les plus courantes :", ""),
  @("The following 'patient id' are atypical cases (one HCAI/CAI with early HAI but no overlap):", "Les « identifiants de patient » suivants sont des cas atypiques (un HCAI/CAI avec HAI précoce mais sans chevauchement) :", ""),
  @("The following 'patient id' are problem case (overlapping specimen collection windows):", "Les 'patient id' suivant sont des cas problématiques (fenêtre de collection de spécimens qui se chevauchent)", ""),
  @("The REDCap dataset is empty/in wrong format. Please contact ACORN support.", "Le jeu de données REDCap est vide/au mauvais format. Veuillez contacter l'assistance ACORN.", ""),
  @("The REDCap dataset is in the right format.", "Le jeu de données REDCap est au bon format.", ""),
  @("There are D28 follow-up done before the expected D28 date.", "Il y a des suivis à J28 effectués avant la date prévue à J28.", ""),
  @("There are multiple F02 with identical ACORN ID, admission date, and episode enrolment date.", "Il existe plusieurs F02 avec un ID ACORN, une date d'admission et une date d'enrôlement identiques.", ""),
  @("There are no atypical case (one HCAI/CAI with early HAI but no overlap).", "Il n'y a pas de cas atypique (un HCAI/CAI avec HAI précoce mais pas de chevauchement).", ""),
  @("There are no D28 follow-up done before the expected D28 date.", "Il n'y a pas de suivi à J28 effectué avant la date prévue à J28.", ""),
  @("There are no isolate with valid AST results. Please contact ACORN support.", "Il n'y a pas d'isolat avec des résultats d'AST valides.", ""),
  @("There are no multiple F02 with identical ACORN ID, admission date, and episode enrolment date.", "Il n'existe pas de F02 avec un ID ACORN, une date d'admission et une date d'enrôlement identiques.", ""),
  @("There are no problem case (overlapping specimen collection windows)", "Il n'y a pas de cas problématique (fenêtre de collection de spécimens qui se chevauchent)", ""),
  @("There are rows for which 'specdate' are after today.", "Il se trouve des enregistrements avec des 'sepcdate' dans le futur.", ""),
  @("There are rows with missing 'orgname'.", "Il existe des enregistrements avec des 'orgname' manquant.", ""),
  @("There are rows with missing 'patid'.", "Il existe des enregistrements avec des 'patid' manquant.", ""),
  @("There are rows with missing 'specdate'.", "Il existe des enregistrements avec des 'specdate' manquant.", ""),
  @("There are rows with missing 'specgroup'.", "Il existe des enregistrements avec des 'specgroup' manquant.", ""),
  @("There are rows with missing 'specid'.", "Il existe des enregistrements avec des 'specid' manquant.", ""),
  @("There is a critical issue with clinical data. The issue should be fixed in REDCap.", "Il y a une erreur critiqe avec les données cliniques. Cette erreur doit être réparée dans REDCap.", ""),
  @("There is no data to display for this organism.", "Cette organisme n'a pas de données à visualiser", ""),
  @("There is no HAI survey data", "Il n'y a pas de données d'enquête HAI", ""),
  @("Trying to save .acorn file on server.", "Tentative de sauvegarde du fichier .acorn sur serveur.", ""),
  @("Updated Charlson Comorbidity Index (uCCI)", "Updated Charlson Comorbidity Index (uCCI)", ""),
  @("User", "Utilisateur", ""),
  @("Variables in Table:", "Variables dans la Table:", ""),
  @("Ward Occupancy Rates", "Taux d'occupation des services", ""),
  @("We couldn't download the lab codes file. Please contact ACORN support.", "Le fichier des codes de lab ne peut pas être téléchargé. Merci de contacter l'équipe ACORN.", ""),
  @("We couldn't download the lab data dictionary. Please contact ACORN support", "Le dictionnaire des données de lab ne peut pas être téléchargé. Merci de contacter l'équipe ACORN.", ""),
  @("Welcome", "Bienvenue", ""),
  @("What do you want to do?", "Comment voulez-vous procéder ?", ""),
  @("With Microbiology", "Avec Microbiologie", ""),
  @("Wrong connection credentials.", "Problème avec la connection internet.", ""),
  @("You are running ACORN dashboard", "Vous utilisez le tableau de bord ACORN", ""),
  @("You can check here if it's the latest production release.", "Vous pouvez vérifier ici s'il s'agit de la dernière version de production.", ""),
  @("Your ACORN dashboard is up to date", "Votre tableau de bord ACORN est à jour", ""),
  @("Follow us on Twitter", "Suivez-nous sur Twitter", ""),
  @("Records in Lab data and BSI forms:", "Enregistrements dans les données de laboratoire et les formulaires BSI:", ""),
  @("ACORN data is not of the right format. Only data generated with v2.1 (or later versions) is compatible.", "TBT", "deleted"),
  @("ACORN data is not of the right format. Only data generated with v2.1 is compatible.", "TBT", "deleted"),
  @("HAI point prevalence by", "Point prévalence HAI par", "deleted")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $val = $row[$j]
        if ($val -ne "") {
            $ws.Cells.Item($i + 1, $j + 1).Value = $val
        }
    }
}

Write-Host "Applied translation updates"